$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0.0001388888888888889
$ws.Range("K2").Value = 953
$ws.Range("L2").Value = 0.001906
